# Handoff report generated: mark the Markdown source as ready for handoff and
# record the newly produced XLF handoff files for each target locale.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276  # OLE (BGR) value of RGB(100,149,237) / #6495ED, the workbook's HyperLink style color

# --- Overview sheet: the Status column reflects the same shared text for the .md row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"

# --- zh-cn locale sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = "Ready for handoff"
$zhcn.Hyperlinks.Add(
    $zhcn.Range("C2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/bc574f683cd1df0c7686ae831f21b604cddbe6fe/6dd68bf3-7360-44d0-b101-438cfcb5d96a.5bd9c4a2ea42f8146e5309f96c66349a843e2cb3.zh-cn.xlf",
    "",
    "",
    "6dd68bf3-7360-44d0-b101-438cfcb5d96a.5bd9c4a2ea42f8146e5309f96c66349a843e2cb3.zh-cn.xlf"
) | Out-Null
$zhcn.Range("C2").Font.Underline = 2
$zhcn.Range("C2").Font.Color = $hyperlinkColor
$zhcn.Range("D2").Value = "2016-01-18 11:24:58"
$zhcn.Range("H2").Value = "Include"

# --- de-de locale sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = "Ready for handoff"
$dede.Hyperlinks.Add(
    $dede.Range("C2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/bc574f683cd1df0c7686ae831f21b604cddbe6fe/6dd68bf3-7360-44d0-b101-438cfcb5d96a.5bd9c4a2ea42f8146e5309f96c66349a843e2cb3.de-de.xlf",
    "",
    "",
    "6dd68bf3-7360-44d0-b101-438cfcb5d96a.5bd9c4a2ea42f8146e5309f96c66349a843e2cb3.de-de.xlf"
) | Out-Null
$dede.Range("C2").Font.Underline = 2
$dede.Range("C2").Font.Color = $hyperlinkColor
$dede.Range("D2").Value = "2016-01-18 11:25:08"
$dede.Range("H2").Value = "Include"
